$d = $word.ActiveDocument

$old = "Unless there are six authors or more give all authors" + [char]0x2019 + " names; do not use " + [char]0x201C + "et al." + [char]0x201D + ". Papers that have not been published, even if they have been submitted for publication, should be cited as " + [char]0x201C + "unpublished" + [char]0x201D + " [4]. Papers that have been accepted for publication should be cited as " + [char]0x201C + "in press" + [char]0x201D + " [5]. Capitalize only the first word in "

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
